$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Conversoes2")

# New torque measurements (comparador de torque) - update source cells;
# the dependent conversion formulas in column B recalc automatically.
$ws.Range("A2").Value = 103.6
$ws.Range("A4").Value = 103.6
$ws.Range("A13").Value = 116
$ws.Range("A15").Value = 116

# Scroll the sheet so row 4 is at the top of the view, keeping the
# existing selection (B15) untouched.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
